$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lu")

# Add category values for existing rows (rows 2-8)
$ws.Range("C2").Value = "Multi-component (informational and operational)"
$ws.Range("C3").Value = "Informational"
$ws.Range("C4").Value = "Multi-component (informational and operational)"
$ws.Range("C5").Value = "Multi-component (informational and operational)"
$ws.Range("C6").Value = "Informational"
$ws.Range("C7").Value = "Informational"
$ws.Range("C8").Value = "Multi-component (informational and operational)"

# Add two new rows
$ws.Range("A9").Value = "plymouth"
$ws.Range("C9").Value = "Informational"

$ws.Range("A10").Value = "dr_julian"
$ws.Range("C10").Value = "Multi-component (informational and operational)"

# Expand the table to include the new rows
$ws.ListObjects.Item("Table1").Resize($ws.Range("A1:C10"))

# Resize column C to fit the new longer category text (matches bestFit width change)
$ws.Columns.Item(3).ColumnWidth = 45.5

# Update selection to match final state
$ws.Range("B7").Select() | Out-Null
